$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.017.36"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.674.10"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.10"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.93%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "21.45"
$c.ClearFormats()
$ws.Range("E9").Value = "  +5.49%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.910.55"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.688.68"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +0.75%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.55%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "66.27"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "27.016.48"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("E18").Value = "  +2.81%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "235.77"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("E24").Value = "  -2.29%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "147.57"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "1.542.72"
$ws.Range("E33").Value = "  +6.73%  "
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("E35").Value = "  +5.00%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("E37").Value = "  +1.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.913"
$c.ClearFormats()
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +2.09%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.05"
$c.ClearFormats()
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("E41").Value = "  +0.05%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "67.72"
$c.ClearFormats()
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  -3.28%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.24"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").Value = "1.817.62"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -0.22%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "90.61"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E50").Value = "  +5.71%  "
$ws.Range("E51").Value = "  +0.30%  "
